# Update Sheets via scheduled runner: refresh Leve profit market data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8568.056
$ws.Range("I51").Value = 8577.333000000001
$ws.Range("J51").Value = 8521.666999999999
$ws.Range("K51").Value = 8577.333000000001
$ws.Range("L51").Value = 8521.666999999999
$ws.Range("M51").Value = -8093.333000000001
$ws.Range("N51").Value = -9489.666999999999
$ws.Range("H64").Value = 71435390
$ws.Range("I64").Value = 6375.5
$ws.Range("J64").Value = 100007000
$ws.Range("K64").Value = 6375.5
$ws.Range("L64").Value = 100007000
$ws.Range("M64").Value = -6127.5
$ws.Range("N64").Value = -100007496
$ws.Range("H67").Value = 71435390
$ws.Range("I67").Value = 6375.5
$ws.Range("J67").Value = 100007000
$ws.Range("K67").Value = 6375.5
$ws.Range("L67").Value = 100007000
$ws.Range("M67").Value = -5517.5
$ws.Range("N67").Value = -100008716
$ws.Range("H107").Value = 360.26315
$ws.Range("I107").Value = 353.6111
$ws.Range("K107").Value = 353.6111
$ws.Range("M107").Value = 1566.3889
$ws.Range("H137").Value = 5568174
$ws.Range("I137").Value = 10000912
$ws.Range("K137").Value = 30002736
$ws.Range("M137").Value = -30000186
$ws.Range("H138").Value = 376031.66
$ws.Range("I138").Value = 2712.1875
$ws.Range("K138").Value = 8136.5625
$ws.Range("M138").Value = -2996.5625
$ws.Range("H141").Value = 5688.278
$ws.Range("I141").Value = 2746.0667
$ws.Range("K141").Value = 8238.2001
$ws.Range("M141").Value = -3058.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3040.8677
$ws.Range("I32").Value = 2486.377
$ws.Range("J32").Value = 7872.857
$ws.Range("K32").Value = 2486.377
$ws.Range("L32").Value = 7872.857
$ws.Range("M32").Value = -2199.377
$ws.Range("N32").Value = -8446.857
$ws.Range("H43").Value = 38042.332
$ws.Range("I43").Value = 23750
$ws.Range("J43").Value = 45188.5
$ws.Range("K43").Value = 23750
$ws.Range("L43").Value = 45188.5
$ws.Range("M43").Value = -23437
$ws.Range("N43").Value = -45814.5
$ws.Range("H45").Value = 3240.5454
$ws.Range("I45").Value = 2624.2222
$ws.Range("K45").Value = 2624.2222
$ws.Range("M45").Value = -2247.2222
$ws.Range("H74").Value = 181610.64
$ws.Range("I74").Value = 279403.4
$ws.Range("K74").Value = 279403.4
$ws.Range("M74").Value = -278529.4
$ws.Range("H77").Value = 181610.64
$ws.Range("I77").Value = 279403.4
$ws.Range("K77").Value = 1397017
$ws.Range("M77").Value = -1392649
$ws.Range("H102").Value = 7493
$ws.Range("I102").Value = 9638
$ws.Range("K102").Value = 9638
$ws.Range("M102").Value = -8016

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 32502870
$ws.Range("I105").Value = 2002412
$ws.Range("K105").Value = 2002412
$ws.Range("M105").Value = -2000665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3697.36
$ws.Range("I31").Value = 2619.9697
$ws.Range("J31").Value = 5788.7646
$ws.Range("K31").Value = 2619.9697
$ws.Range("L31").Value = 5788.7646
$ws.Range("M31").Value = -2324.9697
$ws.Range("N31").Value = -6378.7646
$ws.Range("H34").Value = 3697.36
$ws.Range("I34").Value = 2619.9697
$ws.Range("J34").Value = 5788.7646
$ws.Range("K34").Value = 2619.9697
$ws.Range("L34").Value = 5788.7646
$ws.Range("M34").Value = -2417.9697
$ws.Range("N34").Value = -6192.7646
$ws.Range("H132").Value = 2924.1562
$ws.Range("I132").Value = 2594.7646
$ws.Range("K132").Value = 7784.293799999999
$ws.Range("M132").Value = -5254.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4527.5
$ws.Range("I3").Value = 2764.5
$ws.Range("K3").Value = 8293.5
$ws.Range("M3").Value = -8181.5
$ws.Range("H87").Value = 10156.333
$ws.Range("I87").Value = 8484.5
$ws.Range("K87").Value = 25453.5
$ws.Range("M87").Value = -24205.5
$ws.Range("H90").Value = 10156.333
$ws.Range("I90").Value = 8484.5
$ws.Range("K90").Value = 76360.5
$ws.Range("M90").Value = -70120.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5736.8
$ws.Range("J113").Value = 5921.5
$ws.Range("L113").Value = 5921.5
$ws.Range("N113").Value = -10261.5
$ws.Range("H122").Value = 5210.591
$ws.Range("I122").Value = 4154.2144
$ws.Range("J122").Value = 7059.25
$ws.Range("K122").Value = 12462.6432
$ws.Range("L122").Value = 21177.75
$ws.Range("M122").Value = -10012.6432
$ws.Range("N122").Value = -26077.75
$ws.Range("H126").Value = 9635
$ws.Range("I126").Value = 7839.5
$ws.Range("K126").Value = 23518.5
$ws.Range("M126").Value = -21048.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4067
$ws.Range("I7").Value = 4084
$ws.Range("K7").Value = 4084
$ws.Range("M7").Value = -3972
$ws.Range("H40").Value = 42204.89
$ws.Range("I40").Value = 43443.73
$ws.Range("J40").Value = 9995
$ws.Range("K40").Value = 43443.73
$ws.Range("L40").Value = 9995
$ws.Range("M40").Value = -43307.73
$ws.Range("N40").Value = -10267
$ws.Range("H55").Value = 322.63635
$ws.Range("I55").Value = 185.8
$ws.Range("J55").Value = 436.66666
$ws.Range("K55").Value = 185.8
$ws.Range("L55").Value = 436.66666
$ws.Range("M55").Value = -12.80000000000001
$ws.Range("N55").Value = -782.66666
$ws.Range("H61").Value = 2486.85
$ws.Range("I61").Value = 2444.2942
$ws.Range("K61").Value = 2444.2942
$ws.Range("M61").Value = -2242.2942
$ws.Range("H93").Value = 2424.8572
$ws.Range("I93").Value = 3119.25
$ws.Range("K93").Value = 3119.25
$ws.Range("M93").Value = -1871.25
$ws.Range("H112").Value = 61740.25
$ws.Range("J112").Value = 61740.25
$ws.Range("L112").Value = 61740.25
$ws.Range("N112").Value = -64694.25
$ws.Range("H113").Value = 2486.85
$ws.Range("I113").Value = 2444.2942
$ws.Range("K113").Value = 2444.2942
$ws.Range("M113").Value = -274.2941999999998
$ws.Range("H126").Value = 4067
$ws.Range("I126").Value = 4084
$ws.Range("K126").Value = 12252
$ws.Range("M126").Value = -9782
$ws.Range("H136").Value = 7466.5
$ws.Range("I136").Value = 14400
$ws.Range("J136").Value = 3999.75
$ws.Range("K136").Value = 43200
$ws.Range("L136").Value = 11999.25
$ws.Range("M136").Value = -40650
$ws.Range("N136").Value = -17099.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3777.9473
$ws.Range("I96").Value = 4049.2
$ws.Range("J96").Value = 3476.5557
$ws.Range("K96").Value = 4049.2
$ws.Range("L96").Value = 3476.5557
$ws.Range("M96").Value = -2676.2
$ws.Range("N96").Value = -6222.5557
$ws.Range("H113").Value = 948.6667
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 11365914
$ws.Range("I122").Value = 2319.8125
$ws.Range("K122").Value = 6959.4375
$ws.Range("M122").Value = -4509.4375
$ws.Range("H126").Value = 1081.2
$ws.Range("I126").Value = 1081.2
$ws.Range("K126").Value = 3243.6
$ws.Range("M126").Value = -773.6000000000004
$ws.Range("H132").Value = 5715.591
$ws.Range("I132").Value = 5763.6665
$ws.Range("K132").Value = 17290.9995
$ws.Range("M132").Value = -14760.9995
$ws.Range("H136").Value = 50003930
$ws.Range("J136").Value = 10165.5
$ws.Range("L136").Value = 30496.5
$ws.Range("N136").Value = -35596.5
